# Monday 30 Jul signoff backup
# Adds the Week-5-Friday log entries to the "Daily report" sheet (row 6 col F,
# row 7 cols B & C), matching the format (wrap text, vertical-centered,
# left-aligned for the "Monday..Friday" columns, vertical-centered only for
# the "Weekly Notes"-style columns), grows row 7 to fit the new text, and
# moves the sheet's selection from F6 to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily report")
$ws.Activate()

# --- F6: Friday note for week 5 -------------------------------------------
# Same formatting as D6/E6 (horizontal=left, vertical=center, wrap text).
$ws.Range("D6").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = "Week 5 Friday: Resolved Git LFS issues by uninstalling LFS and restructuring repository. Attempted to implement RT30 for RT60 estimation. Encountered difficulties with irStats function parameter adjustments. Experimented with various approaches but faced persistent calculation issues."

# --- Row 7: week 6, Monday + Tuesday notes ---------------------------------
# B7 matches the D6/E6-style formatting (horizontal=left, vertical=center, wrap).
$ws.Range("D6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = "Investigated RT30 calculation issues. Discovered limitations in energy decay curve fitting. Revised focus to troubleshoot RT60 problems in Unity. Experimented with various Unity/Steam Audio settings to reduce noise floor level. Tested different volume levels, mix ratios, and HRTF settings. Planned comprehensive review of all related components and settings for further troubleshooting."

# C7 matches the B6/C6-style formatting (vertical=center, wrap, default horizontal).
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "WFH. Conducted in-depth review of MATLAB code. Researched related topics to enhance understanding of audio analysis and room acoustics. Prepared for further troubleshooting based on new insights."

# Row 7 grows to fit the new wrapped text.
$ws.Rows.Item(7).RowHeight = 130.5

# --- View: scroll/selection -------------------------------------------------
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C6").Select()
